# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-row observations (date, quality,
# volume, prices, unit, origin, $/kg, kg/unit) among rows 2-20 while the
# leading identity columns (A-C, E-K) stay put. Row 15 is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one "observation" for a given row.
$cols = 4, 12, 13, 14, 15, 16, 17, 18, 19, 20   # D, L, M, N, O, P, Q, R, S, T

# Mapping: new row -> source row (values are read from the ORIGINAL sheet
# state and then written into the destination rows), i.e. a permutation of
# rows 2..20.
$map = @{
    2  = 7
    3  = 5
    4  = 17
    5  = 18
    6  = 13
    7  = 16
    8  = 12
    9  = 20
    10 = 19
    11 = 9
    12 = 11
    13 = 2
    14 = 8
    15 = 15
    16 = 10
    17 = 4
    18 = 14
    19 = 3
    20 = 6
}

# Snapshot the original values for the moved columns of every source row
# before any writes happen (the mapping is a permutation, so rows are both
# read from and written to).
$snapshot = @{}
foreach ($r in 2..20) {
    $rowVals = @{}
    foreach ($c in $cols) {
        # NOTE: `.Value` getter is unreliable on this host (returns a
        # descriptor placeholder); `.Value2` reads the real underlying value.
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
